$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 257, shifting existing rows 257-276 down to 258-277.
$ws.Range("A257").EntireRow.Insert()

# Populate the newly inserted row 257 with the new weekly record.
$ws.Range("A257").Value = 3
$ws.Range("B257").Value = "Femacal de La Calera"
$ws.Range("C257").Value = "Coquimbo"
$ws.Range("D257").Value = 44516
$ws.Range("D257").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E257").Value = 5
$ws.Range("F257").Value = 100112028
$ws.Range("G257").Value = "Sandia"
$ws.Range("H257").Value = "Sin especificar"
$ws.Range("I257").Value = "Primera"
$ws.Range("J257").Value = 180
$ws.Range("K257").Value = 600
$ws.Range("L257").Value = 600
$ws.Range("M257").Value = 600
$ws.Range("N257").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O257").Value = "Perú"
$ws.Range("P257").Value = 600
$ws.Range("Q257").Value = 1
$ws.Range("R257").Value = "Hortaliza"
